# Refresh the cryptos price/volume snapshot (Price column D, Volume(1h) column E).
# Values that look like plain numbers (e.g. "1.00", "600.22") are written with a
# leading apostrophe so Excel stores them as text, matching how this sheet's
# Price column already holds text such as "67.864.57" / "3.791.84".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.864.57'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '3.791.84'
$ws.Range("E3").Value = '  -0.55%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'600.22"
$ws.Range("E5").Value = '  -1.30%  '
$ws.Range("D6").Value = "'163.18"
$ws.Range("E6").Value = '  -2.29%  '
$ws.Range("D7").Value = '3.790.91'
$ws.Range("E7").Value = '  -0.58%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").Value = "'6.66"
$ws.Range("E12").Value = '  +5.71%  '
$ws.Range("E13").Value = '  -3.44%  '
$ws.Range("D14").Value = "'35.16"
$ws.Range("E14").Value = '  -2.58%  '
$ws.Range("D15").Value = '4.428.70'
$ws.Range("E15").Value = '  -0.49%  '
$ws.Range("D16").Value = '3.783.94'
$ws.Range("E16").Value = '  -1.18%  '
$ws.Range("D17").Value = '67.860.35'
$ws.Range("E17").Value = '  +0.10%  '
$ws.Range("D18").Value = "'18.16"
$ws.Range("E18").Value = '  -1.97%  '
$ws.Range("E19").Value = '  +2.08%  '
$ws.Range("D20").Value = "'7.01"
$ws.Range("E20").Value = '  -1.20%  '
$ws.Range("D21").Value = "'458.96"
$ws.Range("E21").Value = '  -0.78%  '
$ws.Range("D22").Value = "'9.44"
$ws.Range("E22").Value = '  -4.55%  '
$ws.Range("E23").Value = '  -1.40%  '
$ws.Range("E24").Value = '  -0.55%  '
$ws.Range("E25").Value = '  -5.61%  '
$ws.Range("E26").Value = '  -1.87%  '
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("E28").Value = '  -0.16%  '
$ws.Range("D29").Value = "'9.89"
$ws.Range("E29").Value = '  -1.24%  '
$ws.Range("D30").Value = '3.938.65'
$ws.Range("E30").Value = '  -0.59%  '
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("E32").Value = '  -2.22%  '
$ws.Range("E33").Value = '  -7.75%  '
$ws.Range("D34").Value = "'28.94"
$ws.Range("E34").Value = '  -2.29%  '
$ws.Range("E35").Value = '  +0.03%  '
$ws.Range("D36").Value = "'8.92"
$ws.Range("E36").Value = '  -1.82%  '
$ws.Range("E37").Value = '  -1.00%  '
$ws.Range("E38").Value = '  +4.84%  '
$ws.Range("D39").Value = "'5.80"
$ws.Range("E39").Value = '  -0.01%  '
$ws.Range("D40").Value = "'0.979"
$ws.Range("E40").Value = '  -1.91%  '
$ws.Range("E41").Value = '  -5.91%  '
$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D44").Value = "'43.61"
$ws.Range("E44").Value = '  +1.07%  '
$ws.Range("D45").Value = "'47.10"
$ws.Range("E45").Value = '  -2.15%  '
$ws.Range("D46").Value = "'151.93"
$ws.Range("E46").Value = '  +2.07%  '
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("E48").Value = '  -0.93%  '
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("E50").Value = '  -0.64%  '
$ws.Range("D51").Value = "'26.43"
$ws.Range("E51").Value = '  -5.63%  '
